$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 3468.8
$ws.Range("I4").Value = 2781.3333
$ws.Range("K4").Value = 2781.3333
$ws.Range("M4").Value = -2667.3333

# Row 17
$ws.Range("H17").Value = 421.34375
$ws.Range("J17").Value = 421.34375
$ws.Range("L17").Value = 1264.03125
$ws.Range("N17").Value = -1600.03125

# Row 61
$ws.Range("H61").Value = 699
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 699
$ws.Range("K61").Value = 0
$ws.Range("L61").ClearContents()
$ws.Range("M61").Value = 2097
$ws.Range("N61").Value = -2441

# Row 64
$ws.Range("H64").Value = 5075.615
$ws.Range("I64").Value = 5236.6
$ws.Range("J64").Value = 4975
$ws.Range("K64").Value = 5236.6
$ws.Range("L64").Value = 4975
$ws.Range("M64").Value = -4988.6
$ws.Range("N64").Value = -5471

# Row 67
$ws.Range("H67").Value = 5075.615
$ws.Range("I67").Value = 5236.6
$ws.Range("J67").Value = 4975
$ws.Range("K67").Value = 5236.6
$ws.Range("L67").Value = 4975
$ws.Range("M67").Value = -4378.6
$ws.Range("N67").Value = -6691

# Row 135
$ws.Range("H135").Value = 23195.834
$ws.Range("I135").Value = 8799
$ws.Range("K135").Value = 79191
$ws.Range("M135").Value = -76656

# Row 138
$ws.Range("H138").Value = 2853.7368
$ws.Range("I138").Value = 1027.1666
$ws.Range("J138").Value = 3196.2188
$ws.Range("K138").Value = 3081.4998
$ws.Range("L138").Value = 9588.6564
$ws.Range("M138").Value = 2058.5002
$ws.Range("N138").Value = -19868.6564

# Row 141
$ws.Range("H141").Value = 6775.4443
$ws.Range("I141").Value = 6711.4287
$ws.Range("K141").Value = 20134.2861
$ws.Range("M141").Value = -14954.2861

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 408.0909
$ws.Range("J5").Value = 939
$ws.Range("L5").Value = 939
$ws.Range("N5").Value = -1163

# Row 45
$ws.Range("H45").Value = 2488
$ws.Range("I45").Value = 2488
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2488
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -2111

# Row 61
$ws.Range("H61").Value = 20882494
$ws.Range("I61").Value = 55558824
$ws.Range("J61").Value = 76697.92999999999
$ws.Range("K61").Value = 55558824
$ws.Range("L61").Value = 76697.92999999999
$ws.Range("M61").Value = -55558612
$ws.Range("N61").Value = -77121.92999999999

# Row 132
$ws.Range("H132").Value = 6403.846
$ws.Range("I132").Value = 3515.9048
$ws.Range("J132").Value = 18533.2
$ws.Range("K132").Value = 10547.7144
$ws.Range("L132").Value = 55599.60000000001
$ws.Range("M132").Value = -8017.714399999999
$ws.Range("N132").Value = -60659.60000000001

# Row 136
$ws.Range("H136").Value = 20882494
$ws.Range("I136").Value = 55558824
$ws.Range("J136").Value = 76697.92999999999
$ws.Range("K136").Value = 166676472
$ws.Range("L136").Value = 230093.79
$ws.Range("M136").Value = -166673922
$ws.Range("N136").Value = -235193.79

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 408.0909
$ws.Range("J4").Value = 939
$ws.Range("L4").Value = 939
$ws.Range("N4").Value = -1169

# Row 82
$ws.Range("H82").Value = 8833.333000000001
$ws.Range("I82").Value = 8833.333000000001
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 8833.333000000001
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -8450.333000000001

# Row 85
$ws.Range("H85").Value = 8833.333000000001
$ws.Range("I85").Value = 8833.333000000001
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 8833.333000000001
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -7507.333000000001

# Row 134
$ws.Range("H134").Value = 41105.152
$ws.Range("I134").Value = 1704.5238
$ws.Range("J134").Value = 206587.8
$ws.Range("K134").Value = 5113.5714
$ws.Range("L134").Value = 619763.3999999999
$ws.Range("M134").Value = -2578.5714
$ws.Range("N134").Value = -624833.3999999999

# Row 135
$ws.Range("H135").Value = 60000
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = 0

# Row 22
$ws.Range("H22").Value = 260.57144
$ws.Range("I22").Value = 256.08334
$ws.Range("K22").Value = 256.08334
$ws.Range("M22").Value = 93.91665999999998

# Row 31
$ws.Range("H31").Value = 806597.3
$ws.Range("J31").Value = 1460402.5
$ws.Range("L31").Value = 1460402.5
$ws.Range("N31").Value = -1460992.5

# Row 34
$ws.Range("H34").Value = 806597.3
$ws.Range("J34").Value = 1460402.5
$ws.Range("L34").Value = 1460402.5
$ws.Range("N34").Value = -1460806.5

# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 0

$ws = $wb.Worksheets.Item("CUL")
# Row 37
$ws.Range("H37").Value = 80000
$ws.Range("J37").Value = 80000
$ws.Range("L37").Value = 240000
$ws.Range("N37").Value = -240224

# Row 39
$ws.Range("H39").Value = 763157
$ws.Range("I39").Value = 526315
$ws.Range("K39").Value = 1578945
$ws.Range("M39").Value = -1578651

# Row 68
$ws.Range("H68").Value = 3632.3333
$ws.Range("I68").Value = 2897
$ws.Range("K68").Value = 8691
$ws.Range("M68").Value = -7880

# Row 71
$ws.Range("H71").Value = 3632.3333
$ws.Range("I71").Value = 2897
$ws.Range("K71").Value = 26073
$ws.Range("M71").Value = -22017

# Row 110
$ws.Range("H110").Value = 14838.333
$ws.Range("J110").Value = 14838.333
$ws.Range("L110").Value = 44514.999
$ws.Range("N110").Value = -52694.999

# Row 112
$ws.Range("H112").Value = 6659.8
$ws.Range("I112").Value = 3600.25
$ws.Range("K112").Value = 10800.75
$ws.Range("M112").Value = -9692.75

# Row 114
$ws.Range("H114").Value = 1808.909
$ws.Range("I114").Value = 1689.7
$ws.Range("J114").Value = 3001
$ws.Range("K114").Value = 5069.1
$ws.Range("L114").Value = 9003
$ws.Range("M114").Value = -1815.1
$ws.Range("N114").Value = -15511

# Row 123
$ws.Range("H123").Value = 57933.61
$ws.Range("I123").Value = 3935.6667
$ws.Range("J123").Value = 68733.2
$ws.Range("K123").Value = 11807.0001
$ws.Range("L123").Value = 206199.6
$ws.Range("M123").Value = -9357.000100000001
$ws.Range("N123").Value = -211099.6

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3766.6667
$ws.Range("J80").Value = 4500
$ws.Range("L80").Value = 4500
$ws.Range("N80").Value = -6496

# Row 83
$ws.Range("H83").Value = 3766.6667
$ws.Range("J83").Value = 4500
$ws.Range("L83").Value = 22500
$ws.Range("N83").Value = -32484

# Row 132
$ws.Range("H132").Value = 90912320
$ws.Range("I132").Value = 125002690
$ws.Range("J132").Value = 4666.6665
$ws.Range("K132").Value = 375008070
$ws.Range("L132").Value = 13999.9995
$ws.Range("M132").Value = -375005540
$ws.Range("N132").Value = -19059.9995

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 3509.76
$ws.Range("I46").Value = 3291.389
$ws.Range("K46").Value = 3291.389
$ws.Range("M46").Value = -3103.389

$ws = $wb.Worksheets.Item("WVR")
# Row 92
$ws.Range("H92").Value = 65395
$ws.Range("J92").Value = 65395
$ws.Range("L92").Value = 65395
$ws.Range("N92").Value = -70387

# Row 113
$ws.Range("H113").Value = 906.6316
$ws.Range("I113").Value = 795
$ws.Range("J113").Value = 1148.5
$ws.Range("K113").Value = 2385
$ws.Range("L113").Value = 3445.5
$ws.Range("M113").Value = -215
$ws.Range("N113").Value = -7785.5

# Row 132
$ws.Range("H132").Value = 1305.5161
$ws.Range("I132").Value = 1377.9259
$ws.Range("J132").Value = 816.75
$ws.Range("K132").Value = 4133.7777
$ws.Range("L132").Value = 2450.25
$ws.Range("M132").Value = -1603.7777
$ws.Range("N132").Value = -7510.25
